$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O1").Value = 0.57641252371265828
$ws.Range("BO1").Value = 0.93288534876520612
$ws.Range("BP1").Value = 0.62259649206260503
$ws.Range("C2").Value = 0.94238156887147295
$ws.Range("BA2").Value = 0.82537428123728684
$ws.Range("A3").Value = 0.90857776744956809
$ws.Range("C4").Value = 0.59564866301991515
$ws.Range("C5").Value = 0.92363790892963282
$ws.Range("G5").Value = 0.80302006553299365
$ws.Range("D6").Value = 0.79123943850355305
$ws.Range("E6").Value = 0.99696165797338665
$ws.Range("H6").Value = 0.82511132825207745
$ws.Range("F7").Value = 0.98411217684819685
$ws.Range("AK7").Value = 0.74084123466414242
$ws.Range("AB8").Value = 0.87643143747697461
$ws.Range("BM8").Value = 0.86075471907173673
$ws.Range("G9").Value = 0.76598992227000195
$ws.Range("H10").Value = 0.71896988369134252
$ws.Range("I10").Value = 0.90237793757470097
$ws.Range("L10").Value = 0.95763051443042213
$ws.Range("I12").Value = 0.68311202186998088
$ws.Range("AC12").Value = 0.81651812125869283
$ws.Range("AZ12").Value = 0.95373931898906439
$ws.Range("K13").Value = 0.84650860616132095
$ws.Range("M14").Value = 0.6037777638085795
$ws.Range("O14").Value = 0.977523975889241
$ws.Range("P14").Value = 0.88839759665890594
$ws.Range("M15").Value = 0.76471343332236463
$ws.Range("P15").Value = 0.93387874026421835
$ws.Range("BO15").Value = 0.86387936815792998
$ws.Range("R17").Value = 0.96797069199275576
$ws.Range("S17").Value = 0.74784421934768019
$ws.Range("T18").Value = 0.95932445662095667
$ws.Range("BC18").Value = 0.97765553253272763
$ws.Range("D19").Value = 0.97831523863977532
$ws.Range("H19").Value = 0.95174943009481949
$ws.Range("K20").Value = 0.95498377151947489
$ws.Range("U20").Value = 0.89765911058623948
$ws.Range("V20").Value = 0.82987138407112093
$ws.Range("BG20").Value = 0.84888054686749692
$ws.Range("S21").Value = 0.72607227614500136
$ws.Range("BI21").Value = 0.92367531385165713
$ws.Range("U22").Value = 0.72311986462848243
$ws.Range("W22").Value = 0.86495942426175643
$ws.Range("Y22").Value = 0.83911120180870913
$ws.Range("U23").Value = 0.8310962431291854
$ws.Range("X23").Value = 0.92566355582351145
$ws.Range("T24").Value = 0.79752326907089421
$ws.Range("V24").Value = 0.85137562759578378
$ws.Range("Y24").Value = 0.74771552585574785
$ws.Range("AO25").Value = 0.94238932290495681
$ws.Range("P27").Value = 0.96784950398960912
$ws.Range("Y27").Value = 0.85187950786230926
$ws.Range("AB27").Value = 0.89943472178856365
$ws.Range("BP27").Value = 0.78427524381709168
$ws.Range("Z28").Value = 0.68513186254370462
$ws.Range("AA29").Value = 0.70278718447414645
$ws.Range("AB29").Value = 0.83026872676076091
$ws.Range("AD29").Value = 0.79010521132929801
$ws.Range("AL29").Value = 0.90873099797151324
$ws.Range("AG30").Value = 0.85694789969865748
$ws.Range("AU30").Value = 0.93179365677817205
$ws.Range("AF31").Value = 0.95717108671298656
$ws.Range("AV31").Value = 0.7595440454290433
$ws.Range("AD32").Value = 0.66165226430587532
$ws.Range("BD32").Value = 0.86028170716746599
$ws.Range("N33").Value = 0.57653571387069247
$ws.Range("AF34").Value = 0.68359465084969939
$ws.Range("R35").Value = 0.93544041311332138
$ws.Range("AH35").Value = 0.93315927980153424
$ws.Range("AH36").Value = 0.85975321853527498
$ws.Range("AK36").Value = 0.81430768043148061
$ws.Range("BO36").Value = 0.92951118300559332
$ws.Range("AI37").Value = 0.78709739717464389
$ws.Range("AK38").Value = 0.97791225236266466
$ws.Range("AN38").Value = 0.95894463578791989
$ws.Range("AK39").Value = 0.92789696494562479
$ws.Range("AL39").Value = 0.7107569563345576
$ws.Range("AN39").Value = 0.98473017160078924
$ws.Range("AO39").Value = 0.8575503709081771
$ws.Range("F41").Value = 0.77893168565552084
$ws.Range("AN42").Value = 0.97921125650183138
$ws.Range("AQ42").Value = 0.58535744654139887
$ws.Range("AR42").Value = 0.98484275032008139
$ws.Range("AD43").Value = 0.96725976646402345
$ws.Range("AR43").Value = 0.97273192343797754
$ws.Range("AS43").Value = 0.97226618587862534
$ws.Range("AY44").Value = 0.92196230513509525
$ws.Range("AT45").Value = 0.86612084377235421
$ws.Range("AU45").Value = 0.95090214985537114
$ws.Range("BJ45").Value = 0.89604441386799372
$ws.Range("V46").Value = 0.70557852391333697
$ws.Range("AT48").Value = 0.89953647124635716
$ws.Range("AW48").Value = 0.59170210821615332
$ws.Range("R49").Value = 0.78293920612147772
$ws.Range("Z49").Value = 0.78229538925007391
$ws.Range("AU49").Value = 0.87203526097433981
$ws.Range("AX49").Value = 0.89642852669314155
$ws.Range("AZ50").Value = 0.69460217461840479
$ws.Range("BP50").Value = 0.87330234346463964
$ws.Range("W51").Value = 0.93389752844565654
$ws.Range("BA51").Value = 0.92541080281146004
$ws.Range("J52").Value = 0.97947392091901975
$ws.Range("BB52").Value = 0.85927754259126177
$ws.Range("AZ53").Value = 0.92493357215075633
$ws.Range("BB53").Value = 0.59307613602548281
$ws.Range("BC54").Value = 0.90677663940177866
$ws.Range("BD54").Value = 0.85856788350170377
$ws.Range("BD55").Value = 0.72624677209133304
$ws.Range("BE56").Value = 0.98880141092256923
$ws.Range("P57").Value = 0.80321310712466598
$ws.Range("AE57").Value = 0.73802895309655459
$ws.Range("BC57").Value = 0.92710845165163569
$ws.Range("BF57").Value = 0.66141073483808421
$ws.Range("BD58").Value = 0.82984625027572334
$ws.Range("BE59").Value = 0.53849100677712158
$ws.Range("BF59").Value = 0.99429339646037485
$ws.Range("Q60").Value = 0.99756530384697739
$ws.Range("BF60").Value = 0.84765055736115091
$ws.Range("O61").Value = 0.98237686697124971
$ws.Range("BG61").Value = 0.92411397332389023
$ws.Range("BI62").Value = 0.93999782097618656
$ws.Range("BK62").Value = 0.68251984785398356
$ws.Range("AR63").Value = 0.94345499448551051
$ws.Range("B64").Value = 0.91734552667798186
$ws.Range("AM64").Value = 0.76723603706020294
$ws.Range("BK64").Value = 0.7670875778428452
$ws.Range("BM64").Value = 0.97072504043722274
$ws.Range("AN65").Value = 0.79576085507985284
$ws.Range("BK65").Value = 0.73978258762242644
$ws.Range("K66").Value = 0.79132947725125402
$ws.Range("BM66").Value = 0.99105769783928843
$ws.Range("BO66").Value = 0.56280737080821264
$ws.Range("BH67").Value = 0.97604116962971599
$ws.Range("BM67").Value = 0.75331474270355891
$ws.Range("BP67").Value = 0.89083045333898614
